function TitleCaseConnectors($s) {
    if ($s -eq $null) { return $s }
    $words = $s.Split(" ")
    $connectors = @("de", "del", "la", "el", "los", "y")
    $result = @()
    foreach ($w in $words) {
        if ($connectors -contains $w) {
            $result += ($w.Substring(0,1).ToUpper() + $w.Substring(1))
        } else {
            $result += $w
        }
    }
    return ($result -join " ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns (row 1) to the new snake_case English names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the Spanish connector words ("de", "del", "la", "el", "los", "y")
# inside the state (A) and municipality (B) names, rows 2 through 455.
for ($r = 2; $r -le 455; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null) {
            if ($v.GetType().Name -eq "String") {
                $t = TitleCaseConnectors $v
                $cell.Value = $t
            }
        }
    }
}

# Drop the trailing metadata/footer rows (456-480): sample-size notes, source,
# attribution and date lines that followed the data table.
$ws.Range("A456:D480").EntireRow.Delete()
